$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1546762589928058
$ws.Range("C2").Value = 0.6402877697841727
$ws.Range("J2").Value = 0.01438848920863309
$ws.Range("P2").Value = 0.1115107913669065
$ws.Range("S2").Value = 0.07913669064748201

$ws.Range("B3").Value = 0.00546448087431694
$ws.Range("C3").Value = 0.0273224043715847
$ws.Range("J3").Value = 0.03278688524590164
$ws.Range("P3").Value = 0.7103825136612022
$ws.Range("S3").Value = 0.2240437158469945

$ws.Range("J4").Value = 0.05084745762711865
$ws.Range("P4").Value = 0.711864406779661
$ws.Range("S4").Value = 0.2372881355932203

$ws.Range("B6").Value = 0.06415094339622641
$ws.Range("D6").Value = 0.01509433962264151
$ws.Range("E6").Value = 0.003773584905660377
$ws.Range("F6").Value = 0.08679245283018867
$ws.Range("J6").Value = 0.2415094339622642
$ws.Range("O6").Value = 0.01509433962264151
$ws.Range("Q6").Value = 0.1773584905660377
$ws.Range("R6").Value = 0.01886792452830189
$ws.Range("S6").Value = 0.3773584905660378

$ws.Range("B7").Value = 0.1464646464646465
$ws.Range("D7").Value = 0.02525252525252525
$ws.Range("F7").Value = 0.06060606060606061
$ws.Range("J7").Value = 0.09595959595959595
$ws.Range("O7").Value = 0.01515151515151515
$ws.Range("Q7").Value = 0.1919191919191919
$ws.Range("R7").Value = 0.0505050505050505
$ws.Range("S7").Value = 0.4141414141414141

$ws.Range("B8").Value = 0.08840864440078586
$ws.Range("D8").Value = 0.02946954813359529
$ws.Range("F8").Value = 0.08644400785854617
$ws.Range("J8").Value = 0.07072691552062868
$ws.Range("O8").Value = 0.01178781925343811
$ws.Range("Q8").Value = 0.2023575638506876
$ws.Range("R8").Value = 0.0962671905697446
$ws.Range("S8").Value = 0.4145383104125737

$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("D9").Value = 0.0213903743315508
$ws.Range("F9").Value = 0.0748663101604278
$ws.Range("J9").Value = 0.106951871657754
$ws.Range("O9").Value = 0.0160427807486631
$ws.Range("Q9").Value = 0.267379679144385
$ws.Range("R9").Value = 0.06951871657754011
$ws.Range("S9").Value = 0.3529411764705883

$ws.Range("B10").Value = 0.1006600660066007
$ws.Range("D10").Value = 0.02722772277227723
$ws.Range("E10").Value = 0.00165016501650165
$ws.Range("F10").Value = 0.06930693069306931
$ws.Range("J10").Value = 0.103960396039604
$ws.Range("O10").Value = 0.01402640264026403
$ws.Range("Q10").Value = 0.25
$ws.Range("R10").Value = 0.06930693069306931
$ws.Range("S10").Value = 0.3638613861386139

$ws.Range("G11").Value = 0.1484848484848485
$ws.Range("J11").Value = 0.1212121212121212
$ws.Range("K11").Value = 0.203030303030303
$ws.Range("L11").Value = 0.5151515151515151
$ws.Range("S11").Value = 0.01212121212121212

$ws.Range("G12").Value = 0.7409638554216867
$ws.Range("J12").Value = 0.2168674698795181
$ws.Range("K12").Value = 0.006024096385542169
$ws.Range("L12").Value = 0.006024096385542169
$ws.Range("S12").Value = 0.03012048192771084

$ws.Range("F13").Value = 0.0196078431372549
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2745098039215687
$ws.Range("S13").Value = 0.0392156862745098

$ws.Range("F15").Value = 0.007547169811320755
$ws.Range("H15").Value = 0.1735849056603773
$ws.Range("I15").Value = 0.09433962264150944
$ws.Range("J15").Value = 0.3358490566037736
$ws.Range("K15").Value = 0.06037735849056604
$ws.Range("M15").Value = 0.01132075471698113
$ws.Range("O15").Value = 0.05660377358490566
$ws.Range("S15").Value = 0.260377358490566

$ws.Range("F16").Value = 0.02985074626865672
$ws.Range("H16").Value = 0.1741293532338309
$ws.Range("I16").Value = 0.02487562189054726
$ws.Range("J16").Value = 0.4328358208955224
$ws.Range("K16").Value = 0.154228855721393
$ws.Range("M16").Value = 0.01492537313432836
$ws.Range("O16").Value = 0.04975124378109453
$ws.Range("S16").Value = 0.1194029850746269

$ws.Range("F17").Value = 0.02425373134328358
$ws.Range("H17").Value = 0.2033582089552239
$ws.Range("I17").Value = 0.08582089552238806
$ws.Range("J17").Value = 0.3899253731343283
$ws.Range("K17").Value = 0.09328358208955224
$ws.Range("M17").Value = 0.01865671641791045
$ws.Range("O17").Value = 0.07835820895522388
$ws.Range("S17").Value = 0.1063432835820896

$ws.Range("F18").Value = 0.025
$ws.Range("H18").Value = 0.20625
$ws.Range("I18").Value = 0.09375
$ws.Range("J18").Value = 0.325
$ws.Range("K18").Value = 0.11875
$ws.Range("M18").Value = 0.05625
$ws.Range("O18").Value = 0.06875000000000001
$ws.Range("S18").Value = 0.10625

$ws.Range("F19").Value = 0.02792862684251358
$ws.Range("H19").Value = 0.2234290147401086
$ws.Range("I19").Value = 0.07525213343677269
$ws.Range("J19").Value = 0.3320403413498836
$ws.Range("K19").Value = 0.1140418929402638
$ws.Range("M19").Value = 0.02172226532195501
$ws.Range("O19").Value = 0.09309542280837858
$ws.Range("S19").Value = 0.1124903025601241
